$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template row used to copy cell formats (styles) from, since these new
# survey rows follow the same per-column formatting as the existing data.
$templateRow = 406
# Row 406 has no value (and thus no style) in column M; the nearest prior
# populated M cell (style-wise identical to every other M cell in the
# sheet) is used as the format template for that column instead.
$templateRowM = 385

# Columns used per new row (varies slightly - some optional form fields
# were left blank by the respondent, matching the source diff exactly).
$rowsCols = @{
    407 = @("A","B","C","D","E","F","G","H","I","J","K","L")
    408 = @("A","B","C","D","E","F","G","H","I","J","K","L")
    409 = @("A","B","C","D","E","F","G","H","I","J","K","L")
    410 = @("A","B","C","E","F","G","H","I","J","K","L")
    411 = @("A","B","C","D","E","F","H","I","K","L")
    412 = @("A","B","C","E","F","G","H","I","J","K","L")
}

foreach ($r in 407..412) {
    foreach ($col in $rowsCols[$r]) {
        $ws.Range("$col$templateRow").Copy() | Out-Null
        $ws.Range("$col$r").PasteSpecial(-4122) | Out-Null
    }
}

$ws.Range("M$templateRowM").Copy() | Out-Null
$ws.Range("M409").PasteSpecial(-4122) | Out-Null

# Row 407
$ws.Range("A407").Value = 44246.611578344906
$ws.Range("B407").Value = "Tampere"
$ws.Range("C407").Value = "26-30 v"
$ws.Range("D407").Value = "Mies"
$ws.Range("E407").Value = 7.0
$ws.Range("F407").Value = "Työntekijä / palkollinen"
$ws.Range("G407").Value = 1.0
$ws.Range("H407").Value = "ohjelmistokehittäjä"
$ws.Range("I407").Value = "Pääosin tai kokonaan etätyö"
$ws.Range("J407").Value = 3900.0
$ws.Range("K407").Value = 55000.0
$ws.Range("L407").Value = "Kyllä"

# Row 408
$ws.Range("A408").Value = 44246.61251465278
$ws.Range("B408").Value = "Turku"
$ws.Range("C408").Value = "21-25 v"
$ws.Range("D408").Value = "mies"
$ws.Range("E408").Value = 3.0
$ws.Range("F408").Value = "Työntekijä / palkollinen"
$ws.Range("G408").Value = 1.0
$ws.Range("H408").Value = "Full stack"
$ws.Range("I408").Value = "Pääosin tai kokonaan toimistolla"
$ws.Range("J408").Value = 3200.0
$ws.Range("K408").Value = 40000.0
$ws.Range("L408").Value = "Ei"

# Row 409
$ws.Range("A409").Value = 44246.613068645835
$ws.Range("B409").Value = "Turku"
$ws.Range("C409").Value = "26-30 v"
$ws.Range("D409").Value = "Mies"
$ws.Range("E409").Value = 2.0
$ws.Range("F409").Value = "Työntekijä / palkollinen"
$ws.Range("G409").Value = 1.0
$ws.Range("H409").Value = "Full-stack ohjelmistokehittäjä (junior)"
$ws.Range("I409").Value = "Pääosin tai kokonaan etätyö"
$ws.Range("J409").Value = 2600.0
$ws.Range("K409").Value = 32500.0
$ws.Range("L409").Value = "Ei"
$ws.Range("M409").Value = "ATR Soft"

# Row 410 (no Gender / column D answer)
$ws.Range("A410").Value = 44246.61409989583
$ws.Range("B410").Value = "PK-Seutu (Helsinki, Espoo, Vantaa)"
$ws.Range("C410").Value = "31-35 v"
$ws.Range("E410").Value = 5.0
$ws.Range("F410").Value = "Työntekijä / palkollinen"
$ws.Range("G410").Value = 1.0
$ws.Range("H410").Value = "full-stack"
$ws.Range("I410").Value = "Pääosin tai kokonaan etätyö"
$ws.Range("J410").Value = 2900.0
$ws.Range("K410").Value = 36000.0
$ws.Range("L410").Value = "Ei"

# Row 411 (entrepreneur - no work-time / monthly salary answer)
$ws.Range("A411").Value = 44246.61679134259
$ws.Range("B411").Value = "Viimsi"
$ws.Range("C411").Value = "36-40 v"
$ws.Range("D411").Value = "Mies"
$ws.Range("E411").Value = 20.0
$ws.Range("F411").Value = "Yrittäjä"
$ws.Range("H411").Value = "sysadmin"
$ws.Range("I411").Value = "Pääosin tai kokonaan etätyö"
$ws.Range("K411").Value = 110000.0
$ws.Range("L411").Value = "Kyllä"

# Row 412 (no Gender / column D answer)
$ws.Range("A412").Value = 44246.62107894676
$ws.Range("B412").Value = "Tampere"
$ws.Range("C412").Value = "36-40 v"
$ws.Range("E412").Value = 12.0
$ws.Range("F412").Value = "Työntekijä / palkollinen"
$ws.Range("G412").Value = 1.0
$ws.Range("H412").Value = "Ohjelmistosuunnittelija"
$ws.Range("I412").Value = "Pääosin tai kokonaan toimistolla"
$ws.Range("J412").Value = 3800.0
$ws.Range("K412").Value = 50000.0
$ws.Range("L412").Value = "Ei"
